$d = $word.ActiveDocument

# Update the date heading in the first paragraph
$d.Content.Find.Execute("2025-03-31 Monday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-04-01 Tuesday", 2)

# Update each math-problem cell in the table (row-major order)
$t = $d.Tables.Item(1)
$values = @(
    "9+22=",
    "34+29=",
    "85-37=",
    "56-17=",
    "88-79=",
    "74+18=",
    "96-38=",
    "8+27=",
    "18+5=",
    "91-79=",
    "12-9=",
    "46+15=",
    "34+9=",
    "49+8=",
    "50-45=",
    "7+34=",
    "29+62=",
    "47+44=",
    "16+77=",
    "62-44=",
    "24-5=",
    "83-8=",
    "36-17=",
    "68+25=",
    "31-26=",
    "82-55=",
    "18+58=",
    "21-7=",
    "54-46=",
    "41-25=",
    "95-29=",
    "33-28=",
    "51-18=",
    "34-18=",
    "60-11=",
    "5+49=",
    "53-5=",
    "8+56=",
    "58+27=",
    "60-23=",
    "19+32=",
    "43-15=",
    "60-45=",
    "50-13=",
    "91-7=",
    "19+53=",
    "26+25=",
    "69+14=",
    "92-47=",
    "28-19=",
    "60-25=",
    "86-47=",
    "96-39=",
    "75-59=",
    "18+6=",
    "48-19=",
    "66+8=",
    "35-28=",
    "18+77=",
    "75+17=",
    "92-75=",
    "43-36=",
    "16-7=",
    "49+14=",
    "59+12=",
    "54-47=",
    "7+36=",
    "84-38=",
    "65+17=",
    "63-56=",
    "80-65=",
    "47+7=",
    "50-45=",
    "68+7=",
    "26+49=",
    "18+68=",
    "58+16=",
    "50-2=",
    "73-24=",
    "20-13=",
    "49+22=",
    "18+7=",
    "95-36=",
    "38+8=",
    "16+48=",
    "71-12=",
    "18+13=",
    "36+55=",
    "63-59=",
    "16+68=",
    "29+65=",
    "44+28=",
    "71-36=",
    "34+49=",
    "27+14=",
    "85-77=",
    "35-6=",
    "6+79=",
    "59+4=",
    "61-46="
)

$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx++
    }
}

Write-Host "Done: updated $idx cells"
